$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 1211.4166
$ws.Range("I28").Value = 1386.4
$ws.Range("K28").Value = 1386.4
$ws.Range("M28").Value = -901.4000000000001
# row 100
$ws.Range("H100").Value = 5000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 5000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -6082
# row 132
$ws.Range("H132").Value = 3460
$ws.Range("I132").Value = 3625
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 10875
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -8345
$ws.Range("N132").Value = -13460

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 3131.182
$ws.Range("I45").Value = 3953.625
$ws.Range("J45").Value = 938
$ws.Range("K45").Value = 3953.625
$ws.Range("L45").Value = 938
$ws.Range("M45").Value = -3576.625
$ws.Range("N45").Value = -1692
# row 132
$ws.Range("H132").Value = 2177.9167
$ws.Range("I132").Value = 2229.182
$ws.Range("J132").Value = 1614
$ws.Range("K132").Value = 6687.545999999999
$ws.Range("L132").Value = 4842
$ws.Range("M132").Value = -4157.545999999999
$ws.Range("N132").Value = -9902

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 76
$ws.Range("H76").Value = 109131.25
$ws.Range("J76").Value = 109131.25
$ws.Range("L76").Value = 109131.25
$ws.Range("N76").Value = -109761.25
# row 79
$ws.Range("H79").Value = 109131.25
$ws.Range("J79").Value = 109131.25
$ws.Range("L79").Value = 109131.25
$ws.Range("N79").Value = -111315.25
# row 94
$ws.Range("H94").Value = 911.3333
$ws.Range("I94").Value = 1323.3334
$ws.Range("J94").Value = 499.33334
$ws.Range("K94").Value = 1323.3334
$ws.Range("L94").Value = 499.33334
$ws.Range("M94").Value = -872.3334
$ws.Range("N94").Value = -1401.33334
# row 99
$ws.Range("H99").Value = 1958.125
$ws.Range("I99").Value = 1958.125
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1958.125
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -460.125
$ws.Range("N99").ClearContents()
# row 134
$ws.Range("H134").Value = 2557.5454
$ws.Range("I134").Value = 2557.5454
$ws.Range("K134").Value = 7672.6362
$ws.Range("M134").Value = -5137.6362

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 99
$ws.Range("H99").Value = 2722.1667
$ws.Range("J99").Value = 2750
$ws.Range("L99").Value = 2750
$ws.Range("N99").Value = -5746
# row 126
$ws.Range("H126").Value = 2722.1667
$ws.Range("J126").Value = 2750
$ws.Range("L126").Value = 8250
$ws.Range("N126").Value = -13190
# row 132
$ws.Range("H132").Value = 3790
$ws.Range("I132").Value = 3790
$ws.Range("K132").Value = 11370
$ws.Range("M132").Value = -8840
# row 134
$ws.Range("H134").Value = 3134.3
$ws.Range("I134").Value = 2988.7778
$ws.Range("K134").Value = 8966.3334
$ws.Range("M134").Value = -6431.3334

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 160.1
$ws.Range("I2").Value = 48.6
$ws.Range("J2").Value = 271.6
$ws.Range("K2").Value = 291.6
$ws.Range("L2").Value = 1629.6
$ws.Range("M2").Value = -178.6
$ws.Range("N2").Value = -1855.6
# row 5
$ws.Range("H5").Value = 299.25
$ws.Range("I5").Value = 299.25
$ws.Range("K5").Value = 897.75
$ws.Range("M5").Value = -785.75
# row 12
$ws.Range("H12").Value = 219.8
$ws.Range("I12").Value = 276.2
$ws.Range("J12").Value = 163.4
$ws.Range("K12").Value = 828.5999999999999
$ws.Range("L12").Value = 490.2
$ws.Range("M12").Value = -655.5999999999999
$ws.Range("N12").Value = -836.2
# row 38
$ws.Range("H38").Value = 1874.4445
$ws.Range("I38").Value = 2107.75
$ws.Range("J38").Value = 8
$ws.Range("K38").Value = 6323.25
$ws.Range("L38").Value = 24
$ws.Range("M38").Value = -5976.25
$ws.Range("N38").Value = -718
# row 46
$ws.Range("H46").Value = 1482.4286
$ws.Range("I46").Value = 1400
$ws.Range("K46").Value = 4200
$ws.Range("M46").Value = -4109
# row 135
$ws.Range("H135").Value = 299.25
$ws.Range("I135").Value = 299.25
$ws.Range("K135").Value = 2693.25
$ws.Range("M135").Value = -158.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 11
$ws.Range("H11").Value = 5500000
$ws.Range("I11").Value = 5500000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 5500000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -5499861
$ws.Range("N11").ClearContents()
# row 88
$ws.Range("H88").Value = 50000
$ws.Range("J88").Value = 50000
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50902
# row 91
$ws.Range("H91").Value = 50000
$ws.Range("J91").Value = 50000
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -53120
# row 94
$ws.Range("H94").Value = 31328.455
$ws.Range("J94").Value = 31328.455
$ws.Range("L94").Value = 31328.455
$ws.Range("N94").Value = -32680.455
# row 132
$ws.Range("H132").Value = 3345.2
$ws.Range("I132").Value = 3345.2
$ws.Range("K132").Value = 10035.6
$ws.Range("M132").Value = -7505.599999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 4
$ws.Range("H4").Value = 505000
$ws.Range("I4").Value = 1000000
$ws.Range("K4").Value = 1000000
$ws.Range("M4").Value = -999887
# row 20
$ws.Range("H20").Value = 8720
$ws.Range("I20").Value = 8720
$ws.Range("K20").Value = 8720
$ws.Range("M20").Value = -8494
# row 23
$ws.Range("H23").Value = 450000
$ws.Range("I23").Value = 450000
$ws.Range("K23").Value = 450000
$ws.Range("M23").Value = -449770
# row 28
$ws.Range("H28").Value = 505000
$ws.Range("I28").Value = 1000000
$ws.Range("K28").Value = 1000000
$ws.Range("M28").Value = -999768
# row 37
$ws.Range("H37").Value = 505000
$ws.Range("I37").Value = 1000000
$ws.Range("K37").Value = 1000000
$ws.Range("M37").Value = -999893
# row 68
$ws.Range("H68").Value = 999
$ws.Range("I68").Value = 999
$ws.Range("K68").Value = 999
$ws.Range("M68").Value = -250
# row 71
$ws.Range("H71").Value = 999
$ws.Range("I71").Value = 999
$ws.Range("K71").Value = 4995
$ws.Range("M71").Value = -1251

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 55
$ws.Range("H55").Value = 11259.5
$ws.Range("I55").Value = 5024
$ws.Range("J55").Value = 17495
$ws.Range("K55").Value = 5024
$ws.Range("L55").Value = 17495
$ws.Range("M55").Value = -4747
$ws.Range("N55").Value = -18049
# row 74
$ws.Range("H74").Value = 28499.4
$ws.Range("J74").Value = 28499.4
$ws.Range("L74").Value = 28499.4
$ws.Range("N74").Value = -30371.4
# row 77
$ws.Range("H77").Value = 28499.4
$ws.Range("J77").Value = 28499.4
$ws.Range("L77").Value = 85498.20000000001
$ws.Range("N77").Value = -94858.20000000001
# row 136
$ws.Range("H136").Value = 2985.4348
$ws.Range("I136").Value = 3079.3809
$ws.Range("J136").Value = 1999
$ws.Range("K136").Value = 9238.1427
$ws.Range("L136").Value = 5997
$ws.Range("M136").Value = -6688.1427
$ws.Range("N136").Value = -11097
